$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text from 18/12/2018
#    to 28/02/2019 on the slide master and every slide layout.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "28/02/2019"
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "28/02/2019"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 8: the "Rectangle 10" shape (nested two groups deep, inside
#    "Group 85" > "Group 44") held three centred paragraphs reading
#    "PID Joint" / "Position" / "Controllers". Collapse it to a single
#    centred paragraph reading "C".
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$group85 = $slide8.Shapes.Item(1)
$rect10 = $group85.GroupItems.Item(4)
$rect10.TextFrame.TextRange.Text = "C"

Write-Host "Edit complete"
